$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to track "Номер контракта" / "Имя" / "Фамилия" (buyer
# info) alongside the goods info. This version drops the buyer columns
# and keeps only the goods grid: "Название товара" / "Количество" / "Сумма".
# Deleting columns A:C shifts the goods columns (formerly D:F) left into A:C.
$ws.Range("A1:C1").EntireColumn.Delete()

# After the shift, the last two goods rows ("Урматы" / "Шаурма с солями")
# need to swap places so "Шаурма с солями" is row 4 and "Урматы" is row 5.
$a4 = $ws.Range("A4").Value()
$b4 = $ws.Range("B4").Value()
$c4 = $ws.Range("C4").Value()
$a5 = $ws.Range("A5").Value()
$b5 = $ws.Range("B5").Value()
$c5 = $ws.Range("C5").Value()

$ws.Range("A4").Value = $a5
$ws.Range("B4").Value = $b5
$ws.Range("C4").Value = $c5
$ws.Range("A5").Value = $a4
$ws.Range("B5").Value = $b4
$ws.Range("C5").Value = $c4
